$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column A (the numeric "12"/"8" rank-like values that duplicated
# column F) is removed entirely; every remaining column shifts one to the
# left (B:F -> A:E).
$ws.Columns("A").Delete()

# The header that is now in D1 (was E1, shared string "MODEL_CONDITION")
# loses its underscore.
$ws.Range("D1").Replace("MODEL_CONDITION", "MODELCONDITION")
